$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-85 and append new rows 86-87 with D,J,K,L,M,P values
$data = @(
    @(2,44956,100,2500,2500,2500,2500),
    @(3,44706,90,4700,4700,4700,4700),
    @(4,44747,80,2500,2500,2500,2500),
    @(5,44753,130,2700,3300,2931,2931),
    @(6,44881,100,2500,2500,2500,2500),
    @(7,44767,180,3300,3300,3300,3300),
    @(8,45041,150,4400,4400,4400,4400),
    @(9,44964,108,2500,2500,2500,2500),
    @(10,44707,100,4700,4700,4700,4700),
    @(11,44679,30,5500,5500,5500,5500),
    @(12,44855,30,2500,2500,2500,2500),
    @(13,44447,75,2200,2200,2200,2200),
    @(14,44893,50,2500,2500,2500,2500),
    @(15,44669,60,6250,6250,6250,6250),
    @(16,44756,120,3300,3300,3300,3300),
    @(17,44839,80,2500,2500,2500,2500),
    @(18,44484,40,2200,2200,2200,2200),
    @(19,44453,20,2300,2300,2300,2300),
    @(20,44781,250,2700,2700,2700,2700),
    @(21,44685,60,5000,6000,5333,5333),
    @(22,44203,30,2000,2000,2000,2000),
    @(23,44798,80,2700,2700,2700,2700),
    @(24,44795,120,2700,2700,2700,2700),
    @(25,44809,150,2700,2700,2700,2700),
    @(26,44474,20,1600,1600,1600,1600),
    @(27,44783,90,2700,2700,2700,2700),
    @(28,44837,50,2500,2500,2500,2500),
    @(29,44818,35,2700,2700,2700,2700),
    @(30,44804,100,3300,3300,3300,3300),
    @(31,44931,140,2500,2500,2500,2500),
    @(32,44487,50,2200,2200,2200,2200),
    @(33,44999,14,4400,4400,4400,4400),
    @(34,45043,120,4400,4400,4400,4400),
    @(35,44497,50,2200,2200,2200,2200),
    @(36,44879,200,2500,2500,2500,2500),
    @(37,44959,140,2500,2500,2500,2500),
    @(38,44476,30,2200,2200,2200,2200),
    @(39,44740,50,2500,2500,2500,2500),
    @(40,44776,100,2700,2700,2700,2700),
    @(41,45069,80,4400,4400,4400,4400),
    @(42,44868,80,2500,2500,2500,2500),
    @(43,44496,40,2200,2200,2200,2200),
    @(44,44832,80,2500,2500,2500,2500),
    @(45,44797,200,2700,2700,2700,2700),
    @(46,44452,120,2300,2300,2300,2300),
    @(47,44972,140,3300,3300,3300,3300),
    @(48,45042,80,4400,4400,4400,4400),
    @(49,45086,80,4400,4400,4400,4400),
    @(50,45016,80,4400,4400,4400,4400),
    @(51,44930,90,2500,2500,2500,2500),
    @(52,45015,90,4400,4400,4400,4400),
    @(53,44910,180,2500,2500,2500,2500),
    @(54,45044,50,4400,4400,4400,4400),
    @(55,45054,90,4100,4100,4100,4100),
    @(56,44882,80,2500,2500,2500,2500),
    @(57,44720,100,3600,3600,3600,3600),
    @(58,44755,90,3300,3300,3300,3300),
    @(59,44895,40,2500,2500,2500,2500),
    @(60,44816,90,2700,2700,2700,2700),
    @(61,44965,80,2500,2500,2500,2500),
    @(62,44966,90,3000,3000,3000,3000),
    @(63,44754,50,3300,3300,3300,3300),
    @(64,44483,50,2200,2200,2200,2200),
    @(65,45002,90,4400,4400,4400,4400),
    @(66,44473,140,1600,1600,1600,1600),
    @(67,44769,140,3300,3300,3300,3300),
    @(68,44952,80,3000,3000,3000,3000),
    @(69,44784,180,2700,2700,2700,2700),
    @(70,45085,150,4400,4400,4400,4400),
    @(71,45055,50,4100,4100,4100,4100),
    @(72,44825,30,2700,2700,2700,2700),
    @(73,44970,140,3300,3300,3300,3300),
    @(74,44757,80,3300,3300,3300,3300),
    @(75,44719,80,3600,3600,3600,3600),
    @(76,44677,20,5500,5500,5500,5500),
    @(77,45071,150,4400,4400,4400,4400),
    @(78,45040,120,4400,4400,4400,4400),
    @(79,44771,30,3300,3300,3300,3300),
    @(80,45111,50,3800,3800,3800,3800),
    @(81,44971,50,3300,3300,3300,3300),
    @(82,44811,120,2700,2700,2700,2700),
    @(83,45112,90,3800,3800,3800,3800),
    @(84,45112,50,3800,3800,3800,3800),
    @(85,44749,80,2500,2500,2500,2500),
    @(86,44741,100,2500,2500,2500,2500),
    @(87,45072,50,4400,4400,4400,4400)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
    $ws.Cells.Item($r, 10).Value = $row[2]
    $ws.Cells.Item($r, 11).Value = $row[3]
    $ws.Cells.Item($r, 12).Value = $row[4]
    $ws.Cells.Item($r, 13).Value = $row[5]
    $ws.Cells.Item($r, 16).Value = $row[6]
}

# Fill in the constant columns for the two newly-added rows (86 and 87)
foreach ($r in 86,87) {
    $ws.Cells.Item($r, 1).Value = 10
    $ws.Cells.Item($r, 2).Value = 'Vega Modelo de Temuco'
    $ws.Cells.Item($r, 3).Value = 'La Araucanía'
    $ws.Cells.Item($r, 5).Value = 9
    $ws.Cells.Item($r, 6).Value = 100112042
    $ws.Cells.Item($r, 7).Value = 'Locoto'
    $ws.Cells.Item($r, 8).Value = 'Sin especificar'
    $ws.Cells.Item($r, 9).Value = 'Primera'
    $ws.Cells.Item($r, 14).Value = '$/kilo'
    $ws.Cells.Item($r, 15).Value = 'Región de Arica y Parinacota'
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = 'Hortaliza'
}
